# Velocity group template - fix "Cincinatti" typo to "Cincinnati" in the
# Alias and Group columns for the Cincinnati All Staff / Employees rows.
# (The Name column intentionally keeps the original "Cincinatti" text, per
# the source diff - only Alias (col B) and Group (col C) are corrected.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: Velocity Cincinatti All Staff
$ws.Range("C14").Value = "Cincinnati"
$ws.Range("B14").Value = "velocitycincinnatiallstaff"

# Row 15: Velocity Cincinatti Employees
$ws.Range("C15").Value = "Cincinnati"
$ws.Range("B15").Value = "velocitycincinnatiemployees"

# Update the active selection to B15, matching the saved view state.
$ws.Range("B15").Select()
